$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Insert a new column before D (shifts old D..K to E..L), carrying over the
# formatting of the column immediately to its left (C), matching how Excel
# inserts a blank column.
$ws.Range("D1").EntireColumn.Insert()

# New column header (row 2) for the inserted "R-square" column.
$ws.Range("D2").Value = "R-square"

# New R-square values for uc / fg / ff (rows 3-5).
$ws.Range("D3").Value = 0.8852
$ws.Range("D4").Value = 0.8842
$ws.Range("D5").Value = 0.9139

# The new R-square column's data cells are left-aligned (distinguishing them
# from the rest of the table).
$ws.Range("D3:D5").HorizontalAlignment = -4131

# Widen column C (to fit "R-square"/formatting) and size the new column D.
$ws.Range("C1").EntireColumn.ColumnWidth = 13.833333333333334
$ws.Range("D1").EntireColumn.ColumnWidth = 10.666666666666666

# Match the selection left after editing.
[void]$ws.Range("D5").Select()
